# Delete obsolete printer rows. Deleting whole rows shifts everything
# below them up, which is exactly the effect seen in the target workbook
# (rows 2-6, 52-54 and 56-63 disappear; the rest shift up to close the gaps).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete from the bottom up so earlier row numbers stay valid as rows
# above a deletion shift upward only *after* that deletion executes.
$ws.Range("A56:Q63").EntireRow.Delete()
$ws.Range("A52:Q54").EntireRow.Delete()
$ws.Range("A2:Q6").EntireRow.Delete()
